$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "t-cabs-observation-af" row (row 2). Its display name and
# some of its codes get folded into a new "Atemfrequenz beatmet" row further
# down the table.
$ws.Rows.Item(2).Delete()

# After the deletion, the former row 10 (Beatmungsstunden) is now row 9.
# Insert a brand-new row above it for "t-cabs-observation-atemfrequenz-beatmet"
# and give it the same formatting as the surrounding data rows.
$ws.Rows.Item(9).Insert()
$ws.Range("A8:K8").Copy()
$ws.Range("A9:K9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(9, 1).Value = "t-cabs-observation-atemfrequenz-beatmet"
$ws.Cells.Item(9, 2).Value = "T-CABS Observation AF"
$ws.Cells.Item(9, 3).Value = "Observation Category Codes#vital-signs"
$ws.Cells.Item(9, 4).Value = ""
$ws.Cells.Item(9, 5).Value = "null#9279-1, null#250810003, null#19840-8, null#152490"
$ws.Cells.Item(9, 6).Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Cells.Item(9, 7).Value = "dateTimeĵ, Periodĵ"
$ws.Cells.Item(9, 8).Value = "Quantityĵ"
$ws.Cells.Item(9, 9).Value = "optional"
$ws.Cells.Item(9, 10).Value = ""
$ws.Cells.Item(9, 11).Value = ""
